$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("D3").Value = 44210
$ws.Range("J3").Value = 1450
$ws.Range("K3").Value = 1600
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1650
$ws.Range("P3").Value = 1650

# Row 4 updates
$ws.Range("D4").Value = 44175
$ws.Range("J4").Value = 1400
$ws.Range("K4").Value = 1900
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1950
$ws.Range("P4").Value = 1950
